# Update cryptocurrency price (D) and volume-change (E) cell values
# to reflect the refreshed values from the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellStyle = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("E2").Style = $cellStyle
$cellStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.08"
$ws.Range("D3").Style = $cellStyle
$cellStyle = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E3").Style = $cellStyle
$cellStyle = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E4").Style = $cellStyle
$cellStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.53"
$ws.Range("D5").Style = $cellStyle
$cellStyle = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E5").Style = $cellStyle
$cellStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D6").Style = $cellStyle
$cellStyle = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("E6").Style = $cellStyle
$cellStyle = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E7").Style = $cellStyle
$cellStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.41"
$ws.Range("D8").Style = $cellStyle
$cellStyle = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("E8").Style = $cellStyle
$cellStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("D9").Style = $cellStyle
$cellStyle = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("E9").Style = $cellStyle
$cellStyle = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E10").Style = $cellStyle
$cellStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("D11").Style = $cellStyle
$cellStyle = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("E11").Style = $cellStyle
$cellStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.85"
$ws.Range("D12").Style = $cellStyle
$cellStyle = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("E12").Style = $cellStyle
$cellStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.623.33"
$ws.Range("D13").Style = $cellStyle
$cellStyle = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E13").Style = $cellStyle
$cellStyle = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("E14").Style = $cellStyle
$cellStyle = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("E15").Style = $cellStyle
$cellStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.729.69"
$ws.Range("D16").Style = $cellStyle
$cellStyle = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E16").Style = $cellStyle
$cellStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.56"
$ws.Range("D17").Style = $cellStyle
$cellStyle = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("E17").Style = $cellStyle
$cellStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.70"
$ws.Range("D18").Style = $cellStyle
$cellStyle = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.47%  "
$ws.Range("E18").Style = $cellStyle
$cellStyle = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E19").Style = $cellStyle
$cellStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("D20").Style = $cellStyle
$cellStyle = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("E20").Style = $cellStyle
$cellStyle = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E21").Style = $cellStyle
$cellStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("D22").Style = $cellStyle
$cellStyle = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("E22").Style = $cellStyle
$cellStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.84"
$ws.Range("D23").Style = $cellStyle
$cellStyle = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E23").Style = $cellStyle
$cellStyle = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("E24").Style = $cellStyle
$cellStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.95"
$ws.Range("D25").Style = $cellStyle
$cellStyle = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E25").Style = $cellStyle
$cellStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.14"
$ws.Range("D26").Style = $cellStyle
$cellStyle = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("E26").Style = $cellStyle
$cellStyle = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E27").Style = $cellStyle
$cellStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.19"
$ws.Range("D28").Style = $cellStyle
$cellStyle = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E28").Style = $cellStyle
$cellStyle = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.48%  "
$ws.Range("E29").Style = $cellStyle
$cellStyle = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E30").Style = $cellStyle
$cellStyle = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E31").Style = $cellStyle
$cellStyle = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("E32").Style = $cellStyle
$cellStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.378.94"
$ws.Range("D33").Style = $cellStyle
$cellStyle = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("E33").Style = $cellStyle
$cellStyle = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("E34").Style = $cellStyle
$cellStyle = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("E35").Style = $cellStyle
$cellStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.982"
$ws.Range("D36").Style = $cellStyle
$cellStyle = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E36").Style = $cellStyle
$cellStyle = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E37").Style = $cellStyle
$cellStyle = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E38").Style = $cellStyle
$cellStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.538"
$ws.Range("D39").Style = $cellStyle
$cellStyle = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.63%  "
$ws.Range("E39").Style = $cellStyle
$cellStyle = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("E40").Style = $cellStyle
$cellStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.969"
$ws.Range("D42").Style = $cellStyle
$cellStyle = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E42").Style = $cellStyle
$cellStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.64"
$ws.Range("D43").Style = $cellStyle
$cellStyle = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E43").Style = $cellStyle
$cellStyle = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("E44").Style = $cellStyle
$cellStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.75"
$ws.Range("D45").Style = $cellStyle
$cellStyle = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("E45").Style = $cellStyle
$cellStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.24"
$ws.Range("D46").Style = $cellStyle
$cellStyle = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("E46").Style = $cellStyle
$cellStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.733.52"
$ws.Range("D47").Style = $cellStyle
$cellStyle = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E47").Style = $cellStyle
$cellStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.84"
$ws.Range("D48").Style = $cellStyle
$cellStyle = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("E48").Style = $cellStyle
$cellStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("D49").Style = $cellStyle
$cellStyle = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E49").Style = $cellStyle
$cellStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0968"
$ws.Range("D50").Style = $cellStyle
$cellStyle = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.27%  "
$ws.Range("E50").Style = $cellStyle
$cellStyle = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.17%  "
$ws.Range("E51").Style = $cellStyle
